$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference cell with the default (unstyled) format, used to restore the
# original style index after forcing text entry via a leading apostrophe.
# (Without the apostrophe, Excel would auto-convert numeric-looking strings
# such as "462.31" into actual numbers, which would not match the source data
# where every Coin/Link/Price/Volume cell is stored as text.)
$plainStyle = $ws.Range("B2").Style

$ws.Range('D2').Value = "'53.784.68"
$ws.Range('D2').Style = $plainStyle
$ws.Range('E2').Value = "'  -8.55%  "
$ws.Range('E2').Style = $plainStyle
$ws.Range('D3').Value = "'2.407.59"
$ws.Range('D3').Style = $plainStyle
$ws.Range('E3').Value = "'  -11.70%  "
$ws.Range('E3').Style = $plainStyle
$ws.Range('E4').Value = "'  +0.18%  "
$ws.Range('E4').Style = $plainStyle
$ws.Range('D5').Value = "'462.31"
$ws.Range('D5').Style = $plainStyle
$ws.Range('E5').Value = "'  -8.42%  "
$ws.Range('E5').Style = $plainStyle
$ws.Range('D6').Value = "'131.48"
$ws.Range('D6').Style = $plainStyle
$ws.Range('E6').Value = "'  -7.10%  "
$ws.Range('E6').Style = $plainStyle
$ws.Range('D7').Value = "'0.995"
$ws.Range('D7').Style = $plainStyle
$ws.Range('E7').Value = "'  -0.23%  "
$ws.Range('E7').Style = $plainStyle
$ws.Range('D8').Value = "'0.488"
$ws.Range('D8').Style = $plainStyle
$ws.Range('E8').Value = "'  -8.03%  "
$ws.Range('E8').Style = $plainStyle
$ws.Range('D9').Value = "'2.428.15"
$ws.Range('D9').Style = $plainStyle
$ws.Range('E9').Value = "'  -11.19%  "
$ws.Range('E9').Style = $plainStyle
$ws.Range('D10').Value = "'0.0948"
$ws.Range('D10').Style = $plainStyle
$ws.Range('E10').Value = "'  -9.21%  "
$ws.Range('E10').Style = $plainStyle
$ws.Range('D11').Value = "'5.32"
$ws.Range('D11').Style = $plainStyle
$ws.Range('E11').Value = "'  -12.11%  "
$ws.Range('E11').Style = $plainStyle
$ws.Range('D12').Value = "'0.313"
$ws.Range('D12').Style = $plainStyle
$ws.Range('E12').Value = "'  -10.08%  "
$ws.Range('E12').Style = $plainStyle
$ws.Range('E13').Value = "'  -4.13%  "
$ws.Range('E13').Style = $plainStyle
$ws.Range('D14').Value = "'2.827.96"
$ws.Range('D14').Style = $plainStyle
$ws.Range('E14').Value = "'  -11.82%  "
$ws.Range('E14').Style = $plainStyle
$ws.Range('D15').Value = "'53.703.43"
$ws.Range('D15').Style = $plainStyle
$ws.Range('E15').Value = "'  -8.84%  "
$ws.Range('E15').Style = $plainStyle
$ws.Range('D16').Value = "'0.0000132"
$ws.Range('D16').Style = $plainStyle
$ws.Range('E16').Value = "'  -2.50%  "
$ws.Range('E16').Style = $plainStyle
$ws.Range('D17').Value = "'19.57"
$ws.Range('D17').Style = $plainStyle
$ws.Range('E17').Value = "'  -9.71%  "
$ws.Range('E17').Style = $plainStyle
$ws.Range('D18').Value = "'2.425.06"
$ws.Range('D18').Style = $plainStyle
$ws.Range('E18').Value = "'  -11.08%  "
$ws.Range('E18').Style = $plainStyle
$ws.Range('D19').Value = "'4.17"
$ws.Range('D19').Style = $plainStyle
$ws.Range('E19').Value = "'  -12.42%  "
$ws.Range('E19').Style = $plainStyle
$ws.Range('D20').Value = "'307.15"
$ws.Range('D20').Style = $plainStyle
$ws.Range('E20').Value = "'  -10.77%  "
$ws.Range('E20').Style = $plainStyle
$ws.Range('D21').Value = "'9.39"
$ws.Range('D21').Style = $plainStyle
$ws.Range('E21').Value = "'  -14.73%  "
$ws.Range('E21').Style = $plainStyle
$ws.Range('D22').Value = "'0.994"
$ws.Range('D22').Style = $plainStyle
$ws.Range('E22').Value = "'  -0.32%  "
$ws.Range('E22').Style = $plainStyle
$ws.Range('D23').Value = "'5.69"
$ws.Range('D23').Style = $plainStyle
$ws.Range('E23').Value = "'  +1.43%  "
$ws.Range('E23').Style = $plainStyle
$ws.Range('D24').Value = "'5.33"
$ws.Range('D24').Style = $plainStyle
$ws.Range('E24').Value = "'  -14.61%  "
$ws.Range('E24').Style = $plainStyle
$ws.Range('D25').Value = "'55.96"
$ws.Range('D25').Style = $plainStyle
$ws.Range('E25').Value = "'  -11.37%  "
$ws.Range('E25').Style = $plainStyle
$ws.Range('E26').Value = "'  +0.78%  "
$ws.Range('E26').Style = $plainStyle
$ws.Range('D27').Value = "'0.385"
$ws.Range('D27').Style = $plainStyle
$ws.Range('E27').Value = "'  -9.80%  "
$ws.Range('E27').Style = $plainStyle
$ws.Range('B28').Value = "'Kaspa"
$ws.Range('B28').Style = $plainStyle
$ws.Range('C28').Value = "'https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range('C28').Style = $plainStyle
$ws.Range('D28').Value = "'0.154"
$ws.Range('D28').Style = $plainStyle
$ws.Range('E28').Value = "'  -10.40%  "
$ws.Range('E28').Style = $plainStyle
$ws.Range('B29').Value = "'WrappedeETH"
$ws.Range('B29').Style = $plainStyle
$ws.Range('C29').Value = "'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range('C29').Style = $plainStyle
$ws.Range('D29').Value = "'2.510.04"
$ws.Range('D29').Style = $plainStyle
$ws.Range('E29').Value = "'  -12.08%  "
$ws.Range('E29').Style = $plainStyle
$ws.Range('D30').Value = "'7.07"
$ws.Range('D30').Style = $plainStyle
$ws.Range('E30').Value = "'  -5.75%  "
$ws.Range('E30').Style = $plainStyle
$ws.Range('D31').Value = "'0.996"
$ws.Range('D31').Style = $plainStyle
$ws.Range('E31').Value = "'  -0.21%  "
$ws.Range('E31').Style = $plainStyle
$ws.Range('D32').Value = "'0.0₃0715"
$ws.Range('D32').Style = $plainStyle
$ws.Range('E32').Value = "'  -14.10%  "
$ws.Range('E32').Style = $plainStyle
$ws.Range('D33').Value = "'145.86"
$ws.Range('D33').Style = $plainStyle
$ws.Range('E33').Value = "'  -3.10%  "
$ws.Range('E33').Style = $plainStyle
$ws.Range('D34').Value = "'17.66"
$ws.Range('D34').Style = $plainStyle
$ws.Range('E34').Value = "'  -7.88%  "
$ws.Range('E34').Style = $plainStyle
$ws.Range('E35').Value = "'  -11.37%  "
$ws.Range('E35').Style = $plainStyle
$ws.Range('E36').Value = "'  -8.41%  "
$ws.Range('E36').Style = $plainStyle
$ws.Range('D37').Value = "'3.51"
$ws.Range('D37').Style = $plainStyle
$ws.Range('E37').Value = "'  -16.20%  "
$ws.Range('E37').Style = $plainStyle
$ws.Range('E38').Value = "'  -6.39%  "
$ws.Range('E38').Style = $plainStyle
$ws.Range('D39').Value = "'0.798"
$ws.Range('D39').Style = $plainStyle
$ws.Range('E39').Value = "'  -16.11%  "
$ws.Range('E39').Style = $plainStyle
$ws.Range('D40').Value = "'0.994"
$ws.Range('D40').Style = $plainStyle
$ws.Range('E40').Value = "'  -0.07%  "
$ws.Range('E40').Style = $plainStyle
$ws.Range('D41').Value = "'32.74"
$ws.Range('D41').Style = $plainStyle
$ws.Range('E41').Value = "'  -8.98%  "
$ws.Range('E41').Style = $plainStyle
$ws.Range('D42').Value = "'0.593"
$ws.Range('D42').Style = $plainStyle
$ws.Range('E42').Value = "'  -1.33%  "
$ws.Range('E42').Style = $plainStyle
$ws.Range('D43').Value = "'0.0522"
$ws.Range('D43').Style = $plainStyle
$ws.Range('E43').Value = "'  -6.01%  "
$ws.Range('E43').Style = $plainStyle
$ws.Range('D44').Value = "'3.24"
$ws.Range('D44').Style = $plainStyle
$ws.Range('E44').Value = "'  -8.17%  "
$ws.Range('E44').Style = $plainStyle
$ws.Range('E45').Value = "'  -2.50%  "
$ws.Range('E45').Style = $plainStyle
$ws.Range('E46').Value = "'  -11.14%  "
$ws.Range('E46').Style = $plainStyle
$ws.Range('D47').Value = "'1.926.00"
$ws.Range('D47').Style = $plainStyle
$ws.Range('E47').Value = "'  -11.95%  "
$ws.Range('E47').Style = $plainStyle
$ws.Range('D48').Value = "'0.0869"
$ws.Range('D48').Style = $plainStyle
$ws.Range('E48').Value = "'  -1.87%  "
$ws.Range('E48').Style = $plainStyle
$ws.Range('D49').Value = "'0.0216"
$ws.Range('D49').Style = $plainStyle
$ws.Range('E49').Value = "'  -4.28%  "
$ws.Range('E49').Style = $plainStyle
$ws.Range('D50').Value = "'4.21"
$ws.Range('D50').Style = $plainStyle
$ws.Range('E50').Value = "'  -10.56%  "
$ws.Range('E50').Style = $plainStyle
$ws.Range('D51').Value = "'16.43"
$ws.Range('D51').Style = $plainStyle
$ws.Range('E51').Value = "'  -13.72%  "
$ws.Range('E51').Style = $plainStyle
